$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Register"
$ws2 = $wb.Worksheets.Item(2)   # "NewUserRegister"

# ---------------------------------------------------------------------------
# Sheet 1 ("Register") - row 2 data updates
# ---------------------------------------------------------------------------

# Phone Number: 9876543210 -> 8876543210 (kept as text)
$ws1.Range("D2").NumberFormat = "@"
$ws1.Range("D2").Value = "8876543210"

# FirstName: Vaibhav -> Rgghav
$ws1.Range("A2").Value = "Rgghav"

# Email: vaibhav@example.com -> Raghav11@example.com, with a new mailto hyperlink
$ws1.Range("C2").Value = "Raghav11@example.com"
$ws1.Hyperlinks.Add($ws1.Range("C2"), "mailto:Raghav11@example.com")
$ws1.Range("C2").Style = "Hyperlink"

# Update the stored selection for sheet1 to C2, while keeping sheet2 active
$ws1.Range("C2").Select()
$ws2.Activate()

# ---------------------------------------------------------------------------
# Sheet 2 ("NewUserRegister") - row 2 data updates
# ---------------------------------------------------------------------------

# LastName: Tore -> otte
$ws2.Range("B2").Value = "otte"

# Email: vaibhav11@example.com -> vaibhavotte511@example.com (existing hyperlink/style kept)
$ws2.Range("C2").Value = "vaibhavotte511@example.com"

# Phone Number: 9576543210 (number) -> 7476543215 (as text)
$ws2.Range("D2").NumberFormat = "@"
$ws2.Range("D2").Value = "7476543215"
